$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the Heading1
#    title paragraph (top of the document).
# ---------------------------------------------------------------------------
$titleText = "Play African Simba Free: Review and Top Features | Novomatic"

$titleParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $titleText) {
        $titleParaIndex = $i
        break
    }
}

$titlePara = $d.Paragraphs.Item($titleParaIndex)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item($titleParaIndex + 1)
$metaPara.Range.Style = $d.Styles.Item("Normal")

$boldText = "Meta description"
$restText = ": Read our review of African Simba, a top online slot game by Novomatic featuring 243 paylines, free spins, and a gamble feature. Play free now."

$insertionPoint = $d.Range($metaPara.Range.Start, $metaPara.Range.Start)
$insertionPoint.InsertAfter($boldText + $restText)

$metaPara = $d.Paragraphs.Item($titleParaIndex + 1)
$boldRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start + $boldText.Length)
$boldRange.Font.Bold = 1

# ---------------------------------------------------------------------------
# 2) Remove the duplicated bold "Play African Simba Free..." paragraph that
#    used to sit right before the closing italic meta-description paragraph
#    near the end of the document.
# ---------------------------------------------------------------------------
$dupTitleParaIndex = -1
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $titleText) {
        $dupTitleParaIndex = $i
        break
    }
}

if ($dupTitleParaIndex -gt 0) {
    $dupTitlePara = $d.Paragraphs.Item($dupTitleParaIndex)
    $dupTitlePara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 3) Replace the text of the final (still-italic) paragraph with the new
#    "Prompt: ..." image-generation prompt, keeping the italic formatting
#    and without letting AutoCorrect turn the straight quotes into curly
#    ones.
# ---------------------------------------------------------------------------
$oldMetaImageText = "Read our review of African Simba, a top online slot game by Novomatic featuring 243 paylines, free spins, and a gamble feature. Play free now."

$metaImageParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $oldMetaImageText) {
        $metaImageParaIndex = $i
        break
    }
}

$lastPara = $d.Paragraphs.Item($metaImageParaIndex)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$newPromptText = "Prompt: Create a cartoon-style feature image for Novomatic's " + [char]34 + "African Simba" + [char]34 + " slot game featuring a happy Maya warrior with glasses. The image should be colorful and eye-catching to represent the vibrant African savanna theme of the game. The Maya warrior could be holding a tribal spear or shield, and should be surrounded by the iconic animals of the savanna, such as a lion, elephant, giraffe, or buffalo. Make sure to include the game's title in bold, African-inspired letters to tie in with the theme."

$lastRange.Text = $newPromptText
